$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Determine the last used row in column C (the "Förändrad" date column)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

# All data rows (2..lastRow) in column C currently hold the date serial
# number 46074 (2026-02-21). Bump every one of them forward by one day
# to 46075 (2026-02-22).
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46074) {
        $cell.Value = 46075
    }
}
